$d = $word.ActiveDocument

$pairs = @(
    @("59×53=3127", "76×65=4940"),
    @("80×44=3520", "34×17=578"),
    @("59×12=708",  "41×12=492"),
    @("67×44=2948", "95×96=9120"),
    @("27×24=648",  "48×40=1920"),
    @("29×61=1769", "85×80=6800"),
    @("84×71=5964", "60×22=1320"),
    @("17×18=306",  "21×40=840"),
    @("48×80=3840", "44×36=1584"),
    @("99×92=9108", "82×80=6560"),
    @("98×31=3038", "74×43=3182"),
    @("51×52=2652", "52×14=728"),
    @("11×82=902",  "20×18=360"),
    @("50×31=1550", "20×15=300"),
    @("59×56=3304", "81×87=7047"),
    @("87×15=1305", "18×74=1332"),
    @("47×72=3384", "59×72=4248"),
    @("84×45=3780", "39×35=1365"),
    @("58×18=1044", "70×80=5600"),
    @("92×71=6532", "81×28=2268"),
    @("60×67=4020", "84×86=7224"),
    @("54×18=972",  "75×57=4275"),
    @("76×19=1444", "60×88=5280"),
    @("18×38=684",  "74×50=3700"),
    @("77×88=6776", "53×68=3604")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
